# RMA Receipt Reversal.xlsx - "RMA Details Maintenance Grid" sheet
# Replace RMA-REX3-* test case identifiers with new RMA-7FWQ-* values
# (plus the corresponding generated Salesforce-style record ids).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Column E (RMA#) filled down first for all three rows ...
$ws.Range("E2").Value = "RMA-7FWQ-001"
$ws.Range("E3").Value = "RMA-7FWQ-002"
$ws.Range("E4").Value = "RMA-7FWQ-003"

# ... then the generated child-record id (J) / line number (F) pair for each row.
$ws.Range("J2").Value = "a7s5f000000xKBqAAM"
$ws.Range("F2").Value = "RMA-7FWQ-1-1"

$ws.Range("J3").Value = "a7s5f000000xKBrAAM"
$ws.Range("F3").Value = "RMA-7FWQ-1-2"

$ws.Range("J4").Value = "a7s5f000000xKBsAAM"
$ws.Range("F4").Value = "RMA-7FWQ-1-3"
